$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "hardware"/"sensor" section (rows 8-11) used two merged 2-row blocks
# (BigFin Scientific + a thick-border spacer row, and R-scrape-of-PDF-Forms +
# another spacer row). Un-merge them before removing the now-unneeded spacer
# rows so Excel doesn't choke on partially merged ranges during the delete.
$null = $ws.Range("A8:A9").UnMerge()
$null = $ws.Range("B8:B9").UnMerge()
$null = $ws.Range("C8:C9").UnMerge()
$null = $ws.Range("A10:A11").UnMerge()
$null = $ws.Range("B10:B11").UnMerge()
$null = $ws.Range("C10:C11").UnMerge()

# Remove the two blank spacer rows (11 then 9, highest index first so the
# second delete still targets the right row).
$null = $ws.Rows("11").Delete()
$null = $ws.Rows("9").Delete()

# Restore the editor's last selection/cursor position.
$null = $ws.Range("C16").Select()

Write-Output "Removed spacer rows and merges from the hardware/sensor section"
